$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.45
$ws.Range("H2").Value = 4.2
$ws.Range("I2").Value = 7.5
$ws.Range("L2").Value = 7.5
$ws.Range("AC2").Value = 8.5
$ws.Range("AD2").Value = 8
$ws.Range("AE2").Value = 23
$ws.Range("AH2").Value = 34
$ws.Range("AJ2").Value = 81
$ws.Range("AK2").Value = 51
$ws.Range("AQ2").Value = 23
$ws.Range("BA2").Value = 201
